$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple K-column +1 / -1 tweaks
$ws.Range("K5").Value = 505
$ws.Range("K6").Value = 415
$ws.Range("K8").Value = 447
$ws.Range("K9").Value = 339
$ws.Range("K10").Value = 342
$ws.Range("K11").Value = 525
$ws.Range("K13").Value = 451
$ws.Range("K14").Value = 337
$ws.Range("K15").Value = 343
$ws.Range("K16").Value = 523
$ws.Range("K20").Value = 507
$ws.Range("K21").Value = 414
$ws.Range("K23").Value = 439
$ws.Range("K24").Value = 328
$ws.Range("K25").Value = 517
$ws.Range("K27").Value = 427
$ws.Range("K29").Value = 514
$ws.Range("K31").Value = 429
$ws.Range("K33").Value = 510
$ws.Range("K36").Value = 435
$ws.Range("K39").Value = 512
$ws.Range("K41").Value = 423
$ws.Range("K45").Value = 508
$ws.Range("K47").Value = 443
$ws.Range("K48").Value = 330
$ws.Range("K49").Value = 453
$ws.Range("K50").Value = 318
$ws.Range("K51").Value = 519

# Unit text clean-up
$ws.Range("G9").Value = "t/cap/year"
$ws.Range("G10").Value = "t/cap/year"
$ws.Range("G15").Value = "t/cap/year"

# Row 54: food waste kcal/cap/day -> kcal / year / cap
$ws.Range("G54").Value = "kcal / year / cap"
$ws.Range("H54").Value = 15793.55
$ws.Range("I54").Value = 1180.745
$ws.Range("J54").Value = 13.3759194406921
$ws.Range("K54").Value = 455

# Row 55
$ws.Range("H55").Value = 2305.86
$ws.Range("I55").Value = 1180.745
$ws.Range("J55").Value = 1.952885678110007
$ws.Range("K55").Value = 457

# Row 56
$ws.Range("I56").Value = 1180.745
$ws.Range("J56").Value = 39.65716560307263
$ws.Range("K56").Value = 320

# Row 57
$ws.Range("I57").Value = 1180.745
$ws.Range("J57").Value = 0.01586286624122905
$ws.Range("K57").Value = 321

# Row 58
$ws.Range("I58").Value = 1180.745
$ws.Range("J58").Value = 41.28071683555721
$ws.Range("K58").Value = 322

# Row 59
$ws.Range("I59").Value = 1180.745
$ws.Range("J59").Value = 0.01651499688755828
$ws.Range("K59").Value = 323

# Row 60
$ws.Range("I60").Value = 1180.745
$ws.Range("J60").Value = 0.005547345108384959
$ws.Range("K60").Value = 521

# Row 61
$ws.Range("I61").Value = 1180.745
$ws.Range("J61").Value = 0.04711432188999318
